$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NumberError")

$newValues = @(
    "Wed Feb 15 12:47:04 EST 2023",
    "Wed Feb 15 12:47:12 EST 2023",
    "Wed Feb 15 12:47:21 EST 2023",
    "Wed Feb 15 12:47:30 EST 2023",
    "Wed Feb 15 12:47:40 EST 2023",
    "Wed Feb 15 12:47:49 EST 2023",
    "Wed Feb 15 12:47:58 EST 2023",
    "Wed Feb 15 12:48:07 EST 2023",
    "Wed Feb 15 12:48:16 EST 2023",
    "Wed Feb 15 12:48:25 EST 2023",
    "Wed Feb 15 12:48:34 EST 2023",
    "Wed Feb 15 12:48:43 EST 2023",
    "Wed Feb 15 12:48:53 EST 2023",
    "Wed Feb 15 12:49:02 EST 2023",
    "Wed Feb 15 12:49:10 EST 2023"
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $ws.Range("B$row").Value = $newValues[$i]
}
